# Rotate the "Recorded By" (column G) comma-separated list for every data
# row so that the last name moves to the front of the list, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "backup@backdoor.com, System, system"
# Rows whose value has only a single name are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $lastIdx = $parts.Count - 1
            $rotated = (@($parts[$lastIdx]) + @($parts[0..($lastIdx - 1)])) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
